# October 2021 meal list update ("updated for 27 tarikh")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------------

# Gopal's bazar/deposit amount (row 7) increased 1900 -> 2500
$ws.Range("B7").Value = 2500

# Day 44494 (row 28) and 44495 (row 29): Q (column for one of the members)
# meal count corrected from 2 to 0 for both days
$ws.Range("Q28").Value = 0
$ws.Range("Q29").Value = 0

# Day 44496 (row 30, the 27th / "27 tarikh") previously had no meals recorded;
# fill in the day's meal counts for each member (K..T)
$ws.Range("K30").Value = 2
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 2
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 2
$ws.Range("P30").Value = 2
$ws.Range("Q30").Value = 0
$ws.Range("R30").Value = 2
$ws.Range("S30").Value = 2
$ws.Range("T30").Value = 2

# Recalculate so all dependent formula cells (totals, meal rate, balances,
# etc.) refresh their cached values.
$excel.CalculateFullRebuild()

# --- View state -------------------------------------------------------------
# Move the selection to F42 and scroll the window so row 19 / column H is at
# the top-left of the visible area.
$ws.Range("F42").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 8
